# "new format for excel"
#
# The guild-config sheet used boolean (TRUE/FALSE) cells for the B:G
# columns of rows 3-9. The new format stores these as plain numeric 0/1
# values instead of booleans, so re-write the range with a numeric 0 -
# this drops the boolean cell type (t="b") and writes a numeric cell
# (no explicit type attribute) while keeping the existing cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:G9").Value = 0
